$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.078.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.176.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.62%  '

$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.164.21'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.503'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.85'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.158'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000227'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.700.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.266.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.191.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.113'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '526.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.728'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +19.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.73%  '

$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '550.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0439'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0838'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.126'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.179.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.274'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.82%  '

$ws.Range("D49").Value = '0.0₃0536'
$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.110'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.35%  '
